# Scene 27A edits
# iron_native Word COM-interop script

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

function Insert-ParagraphBefore($anchorText, $newParaText) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $start = $rng.Start
    $rng.InsertParagraphBefore()
    $ip = $d.Range($start, $start)
    $ip.Text = $newParaText
}

# 1. "Teacher (arms_crossed expressionless): So…" -> "Teacher (arms_crossed annoyed): So…"
Replace-Text "Teacher (arms_crossed expressionless): So" "Teacher (arms_crossed annoyed): So"

# 2. "Teacher (neutral confused): …" -> "Teacher: …"
Replace-Text "Teacher (neutral confused): …" "Teacher: …"

# 3. "Teacher (neutral concerned): First of all, ... what you did wrong." ->
#    "Teacher (arms_crossed neutral): First of all, ... how you messed up."
Replace-Text "Teacher (neutral concerned): First of all, whenever you apologize you shouldn" "Teacher (arms_crossed neutral): First of all, whenever you apologize you shouldn"
Replace-Text "in front, since it shows that you have no idea what you did wrong." "in front, since it shows that you have no idea how you messed up."

# 4. Insert new paragraph "Teacher (neutral sigh):" before "She eyes me carefully..."
Insert-ParagraphBefore "She eyes me carefully before slumping" "Teacher (neutral sigh):"

# 5. "Teacher (neutral sincere): I just wanted..." -> "Teacher (neutral neutral): I just wanted..."
Replace-Text "Teacher (neutral sincere): I just wanted to check up on you." "Teacher (neutral neutral): I just wanted to check up on you."

# 6. "Teacher (neutral expressionless): …" -> "Teacher (neutral thinking): …"
Replace-Text "Teacher (neutral expressionless): …" "Teacher (neutral thinking): …"

# 7. "Teacher (neutral sincere): Well, I guess I can't really ask for more from you." ->
#    "Teacher (neutral sigh): Well, I guess I can't really ask for more from you."
Replace-Text "Teacher (neutral sincere): Well, I guess I" "Teacher (neutral sigh): Well, I guess I"

# 8. "Teacher (neutral neutral): Well, I guess I'll let you go..." -> "Teacher (neutral smirk): Well, I guess I'll let you go..."
Replace-Text "Teacher (neutral neutral): Well, I guess I" "Teacher (neutral smirk): Well, I guess I"

# 9. "Teacher (neutral expressionless): You are not." -> "Teacher (neutral neutral): You are not."
Replace-Text "Teacher (neutral expressionless): You are not." "Teacher (neutral neutral): You are not."

# 10. Insert new paragraph "Teacher (neutral gentle):" before "She smiles uncharacteristically softly."
Insert-ParagraphBefore "She smiles uncharacteristically softly." "Teacher (neutral gentle):"

# 11. "Teacher (neutral smiling): Make sure to take care of yourself, alright?" -> "Teacher: Make sure to take care of yourself, alright?"
Replace-Text "Teacher (neutral smiling): Make sure to take care of yourself, alright?" "Teacher: Make sure to take care of yourself, alright?"

# 12. Insert new paragraph "Teacher (exi):" before "I leave the office, ..."
Insert-ParagraphBefore "I leave the office, a little curious" "Teacher (exi):"

# 13. "Well, whatever the reason is Prim's probably been waiting for a while…" ->
#     "Well, whatever the reason is, Prim's probably been waiting for a while…"
#     (adds a comma, and uses a straight apostrophe instead of a curly one -- set
#     directly via Range.Text so Word's smart-quote autocorrect doesn't curl it back)
$rng = $d.Content
$rng.Find.Execute("Well, whatever the reason is Prim’s probably been waiting for a while…", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = "Well, whatever the reason is, Prim's probably been waiting for a while…"

# 14. styles.xml: docDefaults language en_GB -> en
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.LanguageID = "en"
